$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename sheet "GlobalAxes_DDI_PredVsObs" -> "GlobalAxesSettings"
# ------------------------------------------------------------------
$wsAxes = $wb.Worksheets.Item("GlobalAxes_DDI_PredVsObs")
$wsAxes.Name = "GlobalAxesSettings"

# ------------------------------------------------------------------
# 2. GlobalAxesSettings: insert a new first column ("Plot") in front
#    of the existing Type/Dimension/Unit/GridLines/Scaling columns.
# ------------------------------------------------------------------
$wsAxes.Columns.Item(1).Insert()

# Copy header formatting (fill/font/style) from the old column A
# (now shifted to column B) onto the newly inserted column A.
$wsAxes.Range("B1").Copy() | Out-Null
$wsAxes.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsAxes.Range("A1").Value = "Plot"

# Match the original column's width as closely as the engine allows.
$wsAxes.Columns.Item(1).ColumnWidth = 3.33

# ------------------------------------------------------------------
# 3. Projects sheet: rename header "ID" -> "Id"
# ------------------------------------------------------------------
$wsProjects = $wb.Worksheets.Item("Projects")
$wsProjects.Range("A1").Value = "Id"

# ------------------------------------------------------------------
# 4. Restore selections on each sheet, then make "Projects" the
#    active (selected) sheet/tab, matching the saved workbook state.
# ------------------------------------------------------------------
$wsAxes.Activate() | Out-Null
$wsAxes.Range("C6").Select() | Out-Null

$wsProjects.Activate() | Out-Null
$wsProjects.Range("A2").Select() | Out-Null
